$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BILLING INVOICE")

# Update company name and address (shared strings 0 and 1)
$ws.Range("B9").Value = "2GO LOGISTICS"
$ws.Range("A10").Value = "BRGY. BANAY BANAY, KATAPATAN HOMES, CABUYAO LAGUNA"

# Update the first billing date entry
$ws.Range("B16").Value = 45270.30501157408

# Clear row 18: item number, date, shipment description/label and amount
$ws.Range("A18").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("H18").Value = ""

# Clear row 19: SPO label
$ws.Range("D19").Value = ""
